# Update date and recompute answer table for two-digit x two-digit multiplication worksheet
$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-09-27 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-09-28 Sunday", 2) | Out-Null
$d.Content.Find.Execute("83×16=1328", $true, $false, $false, $false, $false, $true, 1, $false, "98×93=9114", 2) | Out-Null
$d.Content.Find.Execute("56×34=1904", $true, $false, $false, $false, $false, $true, 1, $false, "21×14=294", 2) | Out-Null
$d.Content.Find.Execute("27×22=594", $true, $false, $false, $false, $false, $true, 1, $false, "36×14=504", 2) | Out-Null
$d.Content.Find.Execute("16×24=384", $true, $false, $false, $false, $false, $true, 1, $false, "39×43=1677", 2) | Out-Null
$d.Content.Find.Execute("52×68=3536", $true, $false, $false, $false, $false, $true, 1, $false, "82×18=1476", 2) | Out-Null
$d.Content.Find.Execute("26×87=2262", $true, $false, $false, $false, $false, $true, 1, $false, "13×92=1196", 2) | Out-Null
$d.Content.Find.Execute("65×52=3380", $true, $false, $false, $false, $false, $true, 1, $false, "23×79=1817", 2) | Out-Null
$d.Content.Find.Execute("93×12=1116", $true, $false, $false, $false, $false, $true, 1, $false, "59×75=4425", 2) | Out-Null
$d.Content.Find.Execute("13×69=897", $true, $false, $false, $false, $false, $true, 1, $false, "85×33=2805", 2) | Out-Null
$d.Content.Find.Execute("43×44=1892", $true, $false, $false, $false, $false, $true, 1, $false, "15×57=855", 2) | Out-Null
$d.Content.Find.Execute("14×24=336", $true, $false, $false, $false, $false, $true, 1, $false, "39×34=1326", 2) | Out-Null
$d.Content.Find.Execute("44×60=2640", $true, $false, $false, $false, $false, $true, 1, $false, "94×20=1880", 2) | Out-Null
$d.Content.Find.Execute("41×60=2460", $true, $false, $false, $false, $false, $true, 1, $false, "78×36=2808", 2) | Out-Null
$d.Content.Find.Execute("79×25=1975", $true, $false, $false, $false, $false, $true, 1, $false, "45×34=1530", 2) | Out-Null
$d.Content.Find.Execute("26×63=1638", $true, $false, $false, $false, $false, $true, 1, $false, "25×56=1400", 2) | Out-Null
$d.Content.Find.Execute("77×55=4235", $true, $false, $false, $false, $false, $true, 1, $false, "76×11=836", 2) | Out-Null
$d.Content.Find.Execute("85×28=2380", $true, $false, $false, $false, $false, $true, 1, $false, "65×11=715", 2) | Out-Null
$d.Content.Find.Execute("94×63=5922", $true, $false, $false, $false, $false, $true, 1, $false, "72×87=6264", 2) | Out-Null
$d.Content.Find.Execute("64×29=1856", $true, $false, $false, $false, $false, $true, 1, $false, "15×55=825", 2) | Out-Null
$d.Content.Find.Execute("25×98=2450", $true, $false, $false, $false, $false, $true, 1, $false, "30×34=1020", 2) | Out-Null
$d.Content.Find.Execute("22×31=682", $true, $false, $false, $false, $false, $true, 1, $false, "66×61=4026", 2) | Out-Null
$d.Content.Find.Execute("86×60=5160", $true, $false, $false, $false, $false, $true, 1, $false, "43×44=1892", 2) | Out-Null
$d.Content.Find.Execute("89×44=3916", $true, $false, $false, $false, $false, $true, 1, $false, "16×61=976", 2) | Out-Null
$d.Content.Find.Execute("19×93=1767", $true, $false, $false, $false, $false, $true, 1, $false, "43×94=4042", 2) | Out-Null
$d.Content.Find.Execute("41×64=2624", $true, $false, $false, $false, $false, $true, 1, $false, "86×83=7138", 2) | Out-Null
